# Final update to logs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update / reorder the description text in column G for existing rows ---
# (the underlying shared-string table is re-derived automatically from the
#  text we assign here, so we just need the final text per cell)
$ws.Range("G7").Value  = "Setting up work environment, git and modelsim"
$ws.Range("G8").Value  = "Worked on full adder implementation"
$ws.Range("G9").Value  = "Worked on ripple adder implementation"
$ws.Range("G10").Value = "Worked on arithmetic unit implementing Adder, Zero, ExtWord MUX, AltB  and AltBu"
$ws.Range("G11").Value = "Help debug arithmetic unit (errors with sign extension) Output ExtWord was not matching with test bench values"
$ws.Range("G12").Value = "Screenshots of waves of functional simulation"
$ws.Range("G13").Value = "Helping with screenshots of timing simulations"
$ws.Range("G14").Value = "Start working on documentations, screenshot descriptions, etc."
$ws.Range("G15").Value = "Writing up report and proofreading"
$ws.Range("G16").Value = "Cleaning up documentations and finishing up"

# --- Row 16: end time changes ---
$ws.Range("E16").Value = 0.67708333333333337

# --- Row 17: brand new entry gets filled in ---
$ws.Range("C17").Value = 43926
$ws.Range("D17").Value = 0.67708333333333337
$ws.Range("E17").Value = 0.72916666666666663
$ws.Range("G17").Value = "Adding anotations to pdf and submitting"

# --- Update the saved view/selection state ---
# (topLeftCell/zoomScaleNormal are Excel "normal view" window-state hints that
#  this headless runtime's writer does not expose a way to set directly; the
#  selection itself is the reproducible part of the view change.)
$ws.Range("G13").Select()
